$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Time" column recorded values were re-timed: the old 17:26:00 entry
# becomes 17:35:00, and the old (now unused) 17:27:00 entry is removed -
# row 3 now just reuses the corrected 17:35:00 time too.
$ws.Range("A2").Value = "17:35:00"
$ws.Range("A3").Value = "17:35:00"

# Selection moved to D3.
[void]$ws.Range("D3").Select()
